# This workbook contains a table ("쿼리2") backed by a Power Query / Web
# data source. The source diff corresponds to the table having been
# refreshed: the "월별 누적별풍선" (monthly cumulative balloon) counts in
# column C and the "새로고침시간" (refresh time) timestamps in column D
# were updated with newly-fetched values, the active selection moved, and
# the auto-fit width of column D shrank slightly to match the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated data values (column C: 월별 누적별풍선) ---
$ws.Cells.Item(4, 3).Value  = 527343
$ws.Cells.Item(6, 3).Value  = 425840
$ws.Cells.Item(7, 3).Value  = 379068
$ws.Cells.Item(9, 3).Value  = 269330
$ws.Cells.Item(11, 3).Value = 152599
$ws.Cells.Item(12, 3).Value = 79391

# --- Updated refresh timestamp (column D: 새로고침시간) for every data row ---
$newRefreshTime = 46015.520274814815
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 4).Value = $newRefreshTime
}

# --- Column D width shrank slightly (auto-fit re-evaluated on refresh) ---
$ws.Columns("D").ColumnWidth = 18.428571428571427

# --- Active cell/selection moved ---
$ws.Range("G15").Select()
